$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 'Mohammedeltanany@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, servinaz@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg'
$ws.Cells.Item(3, 7).Value = 'Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, servinaz@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg'
$ws.Cells.Item(4, 7).Value = 'hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, manar.montaser@med.asu.edu.eg, gehanadel@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg'
$ws.Cells.Item(5, 7).Value = 'Nadagouda@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Fatmaelhady@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg'
$ws.Cells.Item(6, 7).Value = 'Kerelos.zareef@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, AbeerRagheb@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg'
$ws.Cells.Item(7, 7).Value = 'Shimaa.ashraf@med.asu.edu.eg, aya.saeed@med.asu.edu.eg, Arwa.elnagar@med.asu.edu.eg'
$ws.Cells.Item(8, 7).Value = 'dina.adel@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, Dalia.m.abdelsalam@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg'
$ws.Cells.Item(9, 7).Value = 'Madeha.Saeed@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, merna.said@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg'
$ws.Cells.Item(10, 7).Value = 'mohamed.saleem@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg'
$ws.Cells.Item(11, 7).Value = 'heba.e.mohamed@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg, mariam.e.mohammed@med.asu.edu.eg'
$ws.Cells.Item(12, 7).Value = 'heba.e.mohamed@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg, mariam.e.mohammed@med.asu.edu.eg'
$ws.Cells.Item(13, 7).Value = 'Aya.alaa-eldein@med.asu.edu.eg, ManaratAleslam-Zefan@med.asu.edu.eg, marian.samir@med.asu.edu.eg'
$ws.Cells.Item(15, 7).Value = 'Walaa.h.ghanima@med.asu.edu.eg, Amr-Saeed@med.asu.edu.eg'
$ws.Cells.Item(16, 7).Value = 'abdullah.elagrody@med.asu.edu.eg, eman.samir@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
$ws.Cells.Item(17, 7).Value = 'marinasorial@med.asu.edu.eg, eman.samir@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg'
$ws.Cells.Item(18, 7).Value = 'Mohammedeltanany@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, servinaz@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg'
$ws.Cells.Item(19, 7).Value = 'Veronia.rafat@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, servinaz@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg'
$ws.Cells.Item(20, 7).Value = 'Mohammedeltanany@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, servinaz@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, heba@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg'
$ws.Cells.Item(21, 7).Value = 'Nadagouda@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Fatmaelhady@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg'
$ws.Cells.Item(22, 7).Value = 'Kerelos.zareef@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, AbeerRagheb@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg'
$ws.Cells.Item(23, 7).Value = 'Shimaa.ashraf@med.asu.edu.eg, aya.saeed@med.asu.edu.eg, Arwa.elnagar@med.asu.edu.eg'
$ws.Cells.Item(24, 7).Value = 'dina.adel@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, Dalia.m.abdelsalam@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg'
$ws.Cells.Item(25, 7).Value = 'Madeha.Saeed@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, merna.said@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg'
$ws.Cells.Item(26, 7).Value = 'mohamed.saleem@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg'
$ws.Cells.Item(27, 7).Value = 'heba.e.mohamed@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg, mariam.e.mohammed@med.asu.edu.eg'
$ws.Cells.Item(28, 7).Value = 'heba.e.mohamed@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg, mariam.e.mohammed@med.asu.edu.eg'
$ws.Cells.Item(29, 7).Value = 'Aya.alaa-eldein@med.asu.edu.eg, ManaratAleslam-Zefan@med.asu.edu.eg, marian.samir@med.asu.edu.eg'
$ws.Cells.Item(31, 7).Value = 'Walaa.h.ghanima@med.asu.edu.eg, Amr-Saeed@med.asu.edu.eg'
$ws.Cells.Item(32, 7).Value = 'abdullah.elagrody@med.asu.edu.eg, eman.samir@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
$ws.Cells.Item(33, 7).Value = 'marinasorial@med.asu.edu.eg, eman.samir@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg'
$ws.Cells.Item(34, 7).Value = 'Veronia.rafat@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, servinaz@med.asu.edu.eg, gehanadel@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg'
$ws.Cells.Item(35, 7).Value = 'System, Mohammedeltanany@med.asu.edu.eg, servinaz@med.asu.edu.eg, gehanadel@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, heba@med.asu.edu.eg'
$ws.Cells.Item(36, 7).Value = 'Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, heba@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg'
$ws.Cells.Item(37, 7).Value = 'Nadagouda@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg'
$ws.Cells.Item(40, 7).Value = 'Madeha.Saeed@med.asu.edu.eg, maimustafa@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, merna.said@med.asu.edu.eg, mernamahrous@med.asu.edu.eg, nahed.mosaad@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, sara_atawia@med.asu.edu.eg'
$ws.Cells.Item(41, 7).Value = 'dina.adel@med.asu.edu.eg, Sarah.Abdelmohsen@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, mernamahrous@med.asu.edu.eg, nadia.mostafa@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg'
$ws.Cells.Item(48, 7).Value = 'Remon.Matta@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg'
$ws.Cells.Item(49, 7).Value = 'eman.samir@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg'
$ws.Cells.Item(50, 7).Value = 'Veronia.rafat@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, servinaz@med.asu.edu.eg, gehanadel@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg'
$ws.Cells.Item(51, 7).Value = 'System, Mohammedeltanany@med.asu.edu.eg, servinaz@med.asu.edu.eg, gehanadel@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, heba@med.asu.edu.eg'
$ws.Cells.Item(52, 7).Value = 'Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, heba@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg'
$ws.Cells.Item(53, 7).Value = 'Nadagouda@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg'
$ws.Cells.Item(56, 7).Value = 'Madeha.Saeed@med.asu.edu.eg, maimustafa@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, merna.said@med.asu.edu.eg, mernamahrous@med.asu.edu.eg, nahed.mosaad@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, sara_atawia@med.asu.edu.eg'
$ws.Cells.Item(57, 7).Value = 'dina.adel@med.asu.edu.eg, Sarah.Abdelmohsen@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, mernamahrous@med.asu.edu.eg, nadia.mostafa@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg'
$ws.Cells.Item(64, 7).Value = 'Remon.Matta@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg'
$ws.Cells.Item(65, 7).Value = 'eman.samir@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg'
$ws.Cells.Item(66, 7).Value = 'Veronia.rafat@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, gehanadel@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg'
$ws.Cells.Item(67, 7).Value = 'Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, servinaz@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg'
$ws.Cells.Item(68, 7).Value = 'Veronia.rafat@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm'
$ws.Cells.Item(69, 7).Value = 'Amera.a.saad@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, AbeerRagheb@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg'
$ws.Cells.Item(70, 7).Value = 'Amera.a.saad@med.asu.edu.eg, Nadagouda@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg'
$ws.Cells.Item(71, 7).Value = 'Sara_nabil@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg'
$ws.Cells.Item(72, 7).Value = 'Madeha.Saeed@med.asu.edu.eg, maimustafa@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, merna.said@med.asu.edu.eg, mernamahrous@med.asu.edu.eg, nahed.mosaad@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, sara_atawia@med.asu.edu.eg'
$ws.Cells.Item(73, 7).Value = 'Madeha.Saeed@med.asu.edu.eg, Dalia.m.abdelsalam@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, merna.said@med.asu.edu.eg, nahed.mosaad@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg'
$ws.Cells.Item(74, 7).Value = 'mohamed.saleem@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg'
$ws.Cells.Item(75, 7).Value = 'mona.I.hussein@med.asu.edu.eg, Alaa.A.Mostafa@med.asu.edu.eg'
$ws.Cells.Item(76, 7).Value = 'mona.I.hussein@med.asu.edu.eg, Alaa.A.Mostafa@med.asu.edu.eg'
$ws.Cells.Item(77, 7).Value = 'user@user.com, nourhan.mostafa@med.asu.edu.eg'
$ws.Cells.Item(79, 7).Value = 'Walaa.h.ghanima@med.asu.edu.eg, Amr-Saeed@med.asu.edu.eg'
$ws.Cells.Item(80, 7).Value = 'neveen.nashaat@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, marina_atef@med.asu.edu.eg, Eman_mohamed@med.asu.edu.eg'
$ws.Cells.Item(81, 7).Value = 'marinasorial@med.asu.edu.eg, eman.samir@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg'
$ws.Cells.Item(82, 7).Value = 'Veronia.rafat@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, gehanadel@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg'
$ws.Cells.Item(83, 7).Value = 'Veronia.rafat@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, servinaz@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg'
$ws.Cells.Item(84, 7).Value = 'Mohammedeltanany@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, servinaz@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, heba@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg'
$ws.Cells.Item(85, 7).Value = 'Amera.a.saad@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, AbeerRagheb@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg'
$ws.Cells.Item(86, 7).Value = 'Amera.a.saad@med.asu.edu.eg, Nadagouda@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg'
$ws.Cells.Item(87, 7).Value = 'Sara_nabil@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg'
$ws.Cells.Item(88, 7).Value = 'Madeha.Saeed@med.asu.edu.eg, maimustafa@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, merna.said@med.asu.edu.eg, mernamahrous@med.asu.edu.eg, nahed.mosaad@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, sara_atawia@med.asu.edu.eg'
$ws.Cells.Item(89, 7).Value = 'Madeha.Saeed@med.asu.edu.eg, Dalia.m.abdelsalam@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, merna.said@med.asu.edu.eg, nahed.mosaad@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg'
$ws.Cells.Item(90, 7).Value = 'mohamed.saleem@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg'
$ws.Cells.Item(91, 7).Value = 'mona.I.hussein@med.asu.edu.eg, Alaa.A.Mostafa@med.asu.edu.eg'
$ws.Cells.Item(92, 7).Value = 'mona.I.hussein@med.asu.edu.eg, Alaa.A.Mostafa@med.asu.edu.eg'
$ws.Cells.Item(93, 7).Value = 'user@user.com, nourhan.mostafa@med.asu.edu.eg'
$ws.Cells.Item(95, 7).Value = 'Walaa.h.ghanima@med.asu.edu.eg, Amr-Saeed@med.asu.edu.eg'
$ws.Cells.Item(96, 7).Value = 'neveen.nashaat@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, marina_atef@med.asu.edu.eg, Eman_mohamed@med.asu.edu.eg'
$ws.Cells.Item(98, 7).Value = 'Mohammedeltanany@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, servinaz@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg'
$ws.Cells.Item(100, 7).Value = 'Veronia.rafat@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm'
$ws.Cells.Item(101, 7).Value = 'Nadagouda@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg'
$ws.Cells.Item(102, 7).Value = 'Amera.a.saad@med.asu.edu.eg, Nadagouda@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg'
$ws.Cells.Item(103, 7).Value = 'Shimaa.ashraf@med.asu.edu.eg, aya.saeed@med.asu.edu.eg, Arwa.elnagar@med.asu.edu.eg'
$ws.Cells.Item(104, 7).Value = 'dina.adel@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg'
$ws.Cells.Item(105, 7).Value = 'dina.adel@med.asu.edu.eg, Sarah.Abdelmohsen@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, mernamahrous@med.asu.edu.eg, nadia.mostafa@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg'
$ws.Cells.Item(106, 7).Value = 'mohamed.saleem@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg'
$ws.Cells.Item(112, 7).Value = 'yassmen.ahmed@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, nahla91@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg'
$ws.Cells.Item(114, 7).Value = 'Mohammedeltanany@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, servinaz@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg'
$ws.Cells.Item(116, 7).Value = 'hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, manar.montaser@med.asu.edu.eg, gehanadel@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg'
$ws.Cells.Item(117, 7).Value = 'Nadagouda@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg'
$ws.Cells.Item(118, 7).Value = 'Amera.a.saad@med.asu.edu.eg, Nadagouda@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg'
$ws.Cells.Item(119, 7).Value = 'Shimaa.ashraf@med.asu.edu.eg, aya.saeed@med.asu.edu.eg, Arwa.elnagar@med.asu.edu.eg'
$ws.Cells.Item(120, 7).Value = 'dina.adel@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg'
$ws.Cells.Item(121, 7).Value = 'dina.adel@med.asu.edu.eg, Sarah.Abdelmohsen@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, mernamahrous@med.asu.edu.eg, nadia.mostafa@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg'
$ws.Cells.Item(122, 7).Value = 'mohamed.saleem@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg'
$ws.Cells.Item(128, 7).Value = 'yassmen.ahmed@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, nahla91@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg'
